$wb = $excel.ActiveWorkbook

# --- Add the "is_targeted list" sheet right after "analyte_class list" ---
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$wsIsTargeted = $wb.Worksheets.Add($null, $afterSheet)
$wsIsTargeted.Name = "is_targeted list"
$wsIsTargeted.Range("A1").Value = "'TRUE"
$wsIsTargeted.Range("A1").Style = "Normal"
$wsIsTargeted.Range("A2").Value = "'FALSE"
$wsIsTargeted.Range("A2").Style = "Normal"

# --- Add the "is_technical_replicate list" sheet right after "is_targeted list" ---
$wsIsTechRep = $wb.Worksheets.Add($null, $wsIsTargeted)
$wsIsTechRep.Name = "is_technical_replicate list"
$wsIsTechRep.Range("A1").Value = "'TRUE"
$wsIsTechRep.Range("A1").Style = "Normal"
$wsIsTechRep.Range("A2").Value = "'FALSE"
$wsIsTechRep.Range("A2").Style = "Normal"

# --- Point the "is_targeted" (N) and "is_technical_replicate" (Q) columns on the
#     main sheet at the new lookup lists instead of the old inline TRUE/FALSE list ---
$wsMain = $wb.Worksheets.Item("Export as TSV")

$rngN = $wsMain.Range("N2:N1048576")
$rngN.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$rngN.Validation.ErrorTitle = "Value must come from list"
$rngN.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

$rngQ = $wsMain.Range("Q2:Q1048576")
$rngQ.Validation.Modify(3, 1, 1, "='is_technical_replicate list'!`$A`$1:`$A`$2")
$rngQ.Validation.ErrorTitle = "Value must come from list"
$rngQ.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

# Restore the originally-active sheet/selection.
$wsMain.Activate()
$wsMain.Range("A1").Select() | Out-Null
